$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.62
$ws.Range("H2").Value = 2.46
$ws.Range("I2").Value = 2.8
$ws.Range("L2").Value = 1.29
$ws.Range("T2").Value = 1.54
$ws.Range("U2").Value = 2.12
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 46
$ws.Range("AB2").Value = 16
$ws.Range("AC2").Value = 10.5
$ws.Range("AF2").Value = 24
$ws.Range("AG2").Value = 15.5
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 44
$ws.Range("AJ2").Value = 55
$ws.Range("AK2").Value = 36
$ws.Range("AM2").Value = 85
$ws.Range("AN2").Value = 25
$ws.Range("AO2").Value = 23
$ws.Range("L4").Value = 1.28
$ws.Range("Q4").Value = 1.68
$ws.Range("R4").Value = 1.48
$ws.Range("T4").Value = 1.51
$ws.Range("U4").Value = 2.16
$ws.Range("I5").Value = 1.86
$ws.Range("P5").Value = 2.82
$ws.Range("Q5").Value = 1.48
$ws.Range("R5").Value = 1.75
$ws.Range("T5").Value = 1.53
$ws.Range("U5").Value = 2.62
$ws.Range("V5").Value = 2.16
$ws.Range("X5").Value = 990
$ws.Range("AA5").Value = 980
$ws.Range("AF5").Value = 980
$ws.Range("AI5").Value = 980
$ws.Range("AK5").Value = 980
$ws.Range("AL5").Value = 980
$ws.Range("AN5").Value = 980
$ws.Range("G6").Value = 3.1
$ws.Range("N6").Value = 3.85
$ws.Range("Q6").Value = 1.82
$ws.Range("V6").Value = 1.55
$ws.Range("W6").Value = 1.48
$ws.Range("AE6").Value = 32
$ws.Range("AI6").Value = 44
$ws.Range("AK6").Value = 36
$ws.Range("AL6").Value = 48
$ws.Range("AM6").Value = 100
$ws.Range("F7").Value = 1.78
$ws.Range("G7").Value = 1.97
$ws.Range("H7").Value = 4.1
$ws.Range("I7").Value = 5.1
$ws.Range("J7").Value = 3.8
$ws.Range("K7").Value = 4.9
$ws.Range("L7").Value = 1.01
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 3.95
$ws.Range("O7").Value = 1.27
$ws.Range("P7").Value = 2.04
$ws.Range("Q7").Value = 1.78
$ws.Range("R7").Value = 1.4
$ws.Range("S7").Value = 3
$ws.Range("T7").Value = 1.75
$ws.Range("U7").Value = 2.06
$ws.Range("V7").Value = 1.24
$ws.Range("W7").Value = 2.04
$ws.Range("X7").Value = 21
$ws.Range("Y7").Value = 22
$ws.Range("Z7").Value = 42
$ws.Range("AB7").Value = 11.5
$ws.Range("AC7").Value = 11.5
$ws.Range("AD7").Value = 23
$ws.Range("AE7").Value = 70
$ws.Range("AF7").Value = 14.5
$ws.Range("AG7").Value = 12.5
$ws.Range("AH7").Value = 23
$ws.Range("AI7").Value = 75
$ws.Range("AK7").Value = 23
$ws.Range("AL7").Value = 42
$ws.Range("AN7").Value = 13.5
$ws.Range("AO7").Value = 70
$ws.Range("S8").Value = 4.2
